$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are digit-strings that Excel would otherwise
# auto-convert to a Number (losing a trailing zero or switching to
# scientific notation). Force these specific cells to Text format first
# so the literal string is preserved exactly, matching the source data.
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"

$ws.Range("D2").Value = '59.624.37'
$ws.Range("E2").Value = '  +1.59%  '
$ws.Range("D3").Value = '3.190.84'
$ws.Range("E3").Value = '  +1.19%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '533.59'
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("D6").Value = '144.20'
$ws.Range("E6").Value = '  +2.94%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  +2.23%  '
$ws.Range("D9").Value = '7.33'
$ws.Range("E9").Value = '  -0.18%  '
$ws.Range("E10").Value = '  +1.78%  '
$ws.Range("E11").Value = '  +1.12%  '
$ws.Range("D12").Value = '3.744.00'
$ws.Range("E12").Value = '  +1.27%  '
$ws.Range("D13").Value = '0.138'
$ws.Range("E13").Value = '  -1.21%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '25.77'
$ws.Range("E14").Value = '  -0.42%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '0.0000172'
$ws.Range("E15").Value = '  +0.96%  '
$ws.Range("D16").Value = '59.726.82'
$ws.Range("E16").Value = '  +1.68%  '
$ws.Range("D17").Value = '3.196.06'
$ws.Range("E17").Value = '  +1.33%  '
$ws.Range("D18").Value = '6.22'
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("D19").Value = '13.08'
$ws.Range("E19").Value = '  +0.45%  '
$ws.Range("D20").Value = '8.19'
$ws.Range("E20").Value = '  +0.22%  '
$ws.Range("D21").Value = '366.54'
$ws.Range("E21").Value = '  -1.45%  '
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("E23").Value = '  +1.29%  '
$ws.Range("D24").Value = '69.61'
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("D25").Value = '8.78'
$ws.Range("E25").Value = '  +9.37%  '
$ws.Range("E26").Value = '  +0.64%  '
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.12%  '
$ws.Range("D28").Value = '0.0₃0885'
$ws.Range("E28").Value = '  +1.22%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '22.20'
$ws.Range("E29").Value = '  +1.25%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '1.88'
$ws.Range("E30").Value = '  -0.38%  '
$ws.Range("D31").Value = '6.07'
$ws.Range("E31").Value = '  -1.51%  '
$ws.Range("D32").Value = '5.31'
$ws.Range("E32").Value = '  +2.24%  '
$ws.Range("E33").Value = '  +1.73%  '
$ws.Range("D34").Value = '6.55'
$ws.Range("E34").Value = '  +4.63%  '
$ws.Range("D35").Value = '156.97'
$ws.Range("E35").Value = '  -1.87%  '
$ws.Range("E36").Value = '  -2.35%  '
$ws.Range("D37").Value = '2.786.47'
$ws.Range("E37").Value = '  +4.98%  '
$ws.Range("D38").Value = '25.61'
$ws.Range("E38").Value = '  +1.49%  '
$ws.Range("E39").Value = '  +1.88%  '
$ws.Range("D40").Value = '1.66'
$ws.Range("E40").Value = '  -0.57%  '
$ws.Range("E41").Value = '  +0.84%  '
$ws.Range("D43").Value = '39.16'
$ws.Range("E43").Value = '  +1.06%  '
$ws.Range("D44").Value = '0.710'
$ws.Range("E44").Value = '  +0.28%  '
$ws.Range("D45").Value = '3.233.19'
$ws.Range("E45").Value = '  +1.18%  '
$ws.Range("D46").Value = '0.103'
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("D47").Value = '0.979'
$ws.Range("E47").Value = '  -0.53%  '
$ws.Range("B48").Value = 'SuiNetwork'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D48").Value = '0.802'
$ws.Range("E48").Value = '  +5.35%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").Value = '6.12'
$ws.Range("E49").Value = '  -1.29%  '
$ws.Range("D50").Value = '20.36'
$ws.Range("E50").Value = '  +0.25%  '
$ws.Range("E51").Value = '  +0.03%  '
